# The deck currently ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (the stock Office palette)
#   ppt/theme/theme2.xml -> "Integral"     (the palette actually applied to
#                                            the slide master / the deck's
#                                            live Design)
#
# The authored change swaps the two palettes: the presentation's live theme
# (reached here through SlideMaster.Theme, exactly like a user tweaking
# Design > Variants > Colors > Customize Colors...) becomes the classic
# Office palette, colour-for-colour.
#
# RGB() in the PowerPoint/VBA object model packs R + G*256 + B*65536, so a
# small helper turns the familiar RRGGBB hex from the target theme into the
# long value ThemeColorScheme items expect.

function Hex-ToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Target palette = the "Office Theme" colours, applied in the standard
# msoThemeColor order: Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink,
# FollowedHyperlink.
$officePalette = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officePalette.Length; $i++) {
    $colors.Item($i + 1).RGB = Hex-ToRgbLong $officePalette[$i]
}
